$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "250.73"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.88"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.397"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05648"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.435"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.371"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8177"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9220"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1430"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07498"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03161"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03079"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09325"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.548"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001597"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04699"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005775"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006377"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005013"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001496"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.726"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.180"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3293"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1303"
$ws.Range("E27").Value = "26AAXTokenAABWorstin24h"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04013"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006941"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1068"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002781"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007574"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005557"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000748"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6581"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.2235"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002094"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01007"
